# Bugfixed evaluation and simulated rt_data for components
# Rewrites the data rows (A2:E19) of Sheet1 with corrected forecast values,
# shifting everything down by one row (new first row for 2007 inserted,
# 2024/2025 data added at the bottom) and updating the dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full corrected data set: date_serial, y_0, y_0_forecast, y_1, y_1_forecast
$data = @(
    @(39400, 2007, 7.226520411029069,  2008, $Null),
    @(39765, 2008, 4.268860212333636,  2009, $Null),
    @(40130, 2009, -7.266312015249776, 2010, $Null),
    @(40494, 2010, 6.958243460951929,  2011, 12.21658306395068),
    @(40862, 2011, 9.469137444079934,  2012, 8.079264579851909),
    @(41228, 2012, 3.358206407534947,  2013, 4.701432377325987),
    @(41592, 2013, 0.3081076735359067, 2014, 3.972902167062387),
    @(41957, 2014, 3.901355411819707,  2015, 4.658857392675264),
    @(42321, 2015, 5.331683351557981,  2016, 4.089819750351786),
    @(42689, 2016, 3.254758369308375,  2017, 2.313009565865753),
    @(43053, 2017, 5.246209615995667,  2018, 4.784022165496182),
    @(43418, 2018, 4.86255966374296,   2019, 4.112897401876747),
    @(43783, 2019, 2.764740011159428,  2020, 1.643374185611401),
    @(44159, 2020, -7.260793671746435, 2021, 0.00562230452727519),
    @(44525, 2021, 4.097586525396268,  2022, 3.9116372951149),
    @(44890, 2022, 7.824284864703746,  2023, 2.586378346096296),
    @(45254, 2023, -1.24502235313334,  2024, -1.561801765212567),
    @(45618, 2024, -1.735114423676209, 2025, 2.409056355286521)
)

$startRow = 2
$endRow = $startRow + $data.Length - 1
$lastExistingRow = 18

# Clear out any leftover values beyond the new data range (defensive, in
# case a previous layout had more/less rows).
if ($endRow + 1 -le 1000) {
    $ws.Range("A$($endRow + 1)`:E1000").ClearContents()
}

# The new row(s) beyond what previously existed don't carry any formatting
# yet. Copy the date column's look (bold font, thin border, centered,
# custom date number format) from an existing row so the appended row
# matches the rest of the table.
if ($endRow -gt $lastExistingRow) {
    $formatSource = $ws.Range("A2")
    for ($row = $lastExistingRow + 1; $row -le $endRow; $row++) {
        $formatSource.Copy($ws.Range("A$row"))
    }
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]

    $ws.Cells.Item($row, 2).Value = $vals[1]
    if ($vals[2] -ne $Null) {
        $ws.Cells.Item($row, 3).Value = $vals[2]
    } else {
        $ws.Cells.Item($row, 3).ClearContents()
    }
    $ws.Cells.Item($row, 4).Value = $vals[3]
    if ($vals[4] -ne $Null) {
        $ws.Cells.Item($row, 5).Value = $vals[4]
    } else {
        $ws.Cells.Item($row, 5).ClearContents()
    }
}

